$wb = $excel.ActiveWorkbook

# Mapping of row -> [old, new] value in column F for sheets "展览" and "全部类型"
$updates = @{
    4  = 151
    5  = 7047
    6  = 3979
    10 = 40
    11 = 92
    12 = 56
    13 = 43
    14 = 185
    15 = 605
    16 = 71
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
